$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38
$ws.Range("A38").Value = "2025/26"
$ws.Range("B38").Value = "WSL"

# Row 39
$ws.Range("A39").Value = "2025/26"
$ws.Range("B39").Value = "WSL2"
$ws.Range("C39").Value = "2856p00v9pp48aeyzcso6i32s"
$ws.Range("C38").Value = "221phckhkd7y6rg3uyava3ifo"

# Row 40
$ws.Range("A40").Value = "2025/26"
$ws.Range("C40").Value = "3ielq9pcsvqfftb4q5zjc2dw"
$ws.Range("B40").Value = "A-League Women"

# Row 41
$ws.Range("A41").Value = "2025/26"
$ws.Range("C41").Value = "2bqrpllc5x3it55paifyfa044"
$ws.Range("B41").Value = "Premiere League"

# Row 42
$ws.Range("A42").Value = "2025/26"
$ws.Range("C42").Value = "e98d4oial167tji58n80jkh04"
$ws.Range("B42").Value = "USL Super League"

# Update view: scroll and select next empty cell in column B, matching Excel's
# behaviour after typing data down a column and moving to the next row.
$ws.Range("B43").Select()
$excel.ActiveWindow.ScrollRow = 7
